$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "61.984.77", "1.00", "13.10") are preserved exactly as text,
# matching the inline-string cell type used in the source workbook.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.017.99'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.423.86'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '411.29'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '130.18'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.641'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +8.15%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +8.28%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +14.78%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '43.06'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +61.13%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.13'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +7.76%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.956.32'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.25'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +6.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.424.91'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.24'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +7.41%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +6.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '62.014.40'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '448.33'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +42.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '91.54'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +8.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.18'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.10'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.28'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '33.69'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +13.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.80'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +6.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '4.75'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.72'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +5.50%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '43.04'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.168'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -3.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0501'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '53.93'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +4.81%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +7.69%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.40%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '141.89'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +2.62%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +4.70%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.99'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +7.71%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.72'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.38'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +5.01%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.765.23'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.110.75'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '105.35'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +25.65%  '
